$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 831
$ws.Range("D2").Value = 863

$ws.Range("C3").Value = 1147
$ws.Range("D3").Value = 1203
$ws.Range("E3").Value = 0.49

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").ClearContents()

$ws.Range("C5").Value = 1519
$ws.Range("D5").Value = 1523
$ws.Range("E5").Value = 0.5

$ws.Range("C6").Value = 1940
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 0.5
